$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''243.84'
$ws.Range('G2').Value = '''18'
$ws.Range('G3').Value = '''18'
$ws.Range('D4').Value = '''5.322'
$ws.Range('G4').Value = '''18'
$ws.Range('D5').Value = '''0.05653'
$ws.Range('G5').Value = '''18'
$ws.Range('D6').Value = '''3.376'
$ws.Range('G6').Value = '''18'
$ws.Range('D7').Value = '''6.384'
$ws.Range('G7').Value = '''18'
$ws.Range('D8').Value = '''0.8056'
$ws.Range('G8').Value = '''18'
$ws.Range('D9').Value = '''0.9471'
$ws.Range('G9').Value = '''18'
$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D10').Value = '''0.1421'
$ws.Range('E10').Value = '9WazirXWRX'
$ws.Range('G10').Value = '''18'
$ws.Range('B11').Value = 'MandalaExchangeToken'
$ws.Range('C11').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D11').Value = '''0.07409'
$ws.Range('E11').Value = '10MandalaExchangeTokenMDX'
$ws.Range('G11').Value = '''18'
$ws.Range('B12').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C12').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D12').Value = '''0.03192'
$ws.Range('E12').Value = '11LiechtensteinCryptoassetsExchangeLCX'
$ws.Range('G12').Value = '''18'
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D13').Value = '''0.03055'
$ws.Range('E13').Value = '12BitrueCoinBTR'
$ws.Range('G13').Value = '''18'
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').Value = '''0.09269'
$ws.Range('E14').Value = '13BitMartTokenBMX'
$ws.Range('G14').Value = '''18'
$ws.Range('B15').Value = 'MCDex'
$ws.Range('C15').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('D15').Value = '''3.573'
$ws.Range('E15').Value = '14MCDexMCB'
$ws.Range('G15').Value = '''18'
$ws.Range('B16').Value = 'BitForexToken'
$ws.Range('C16').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D16').Value = '''0.001647'
$ws.Range('E16').Value = '15BitForexTokenBF'
$ws.Range('G16').Value = '''18'
$ws.Range('B17').Value = 'CoinExToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D17').Value = '''0.04697'
$ws.Range('E17').Value = '16CoinExTokenCET'
$ws.Range('G17').Value = '''18'
$ws.Range('B18').Value = 'One'
$ws.Range('C18').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D18').Value = '''0.0005825'
$ws.Range('E18').Value = '17OneONE'
$ws.Range('G18').Value = '''18'
$ws.Range('D19').Value = '''0.006351'
$ws.Range('G19').Value = '''18'
$ws.Range('D20').Value = '''0.004979'
$ws.Range('G20').Value = '''18'
$ws.Range('G21').Value = '''18'
$ws.Range('G22').Value = '''18'
$ws.Range('D23').Value = '''0.0003103'
$ws.Range('G23').Value = '''18'
$ws.Range('D24').Value = '''3.765'
$ws.Range('G24').Value = '''18'
$ws.Range('D25').Value = '''2.098'
$ws.Range('G25').Value = '''18'
$ws.Range('D26').Value = '''0.3271'
$ws.Range('G26').Value = '''18'
$ws.Range('G27').Value = '''18'
$ws.Range('G28').Value = '''18'
$ws.Range('G29').Value = '''18'
$ws.Range('G30').Value = '''18'
$ws.Range('G31').Value = '''18'
$ws.Range('G32').Value = '''18'
$ws.Range('G33').Value = '''18'
$ws.Range('G34').Value = '''18'
$ws.Range('G35').Value = '''18'
$ws.Range('G36').Value = '''18'
$ws.Range('G37').Value = '''18'
$ws.Range('G38').Value = '''18'
$ws.Range('G39').Value = '''18'
$ws.Range('D40').Value = '''0.03944'
$ws.Range('G40').Value = '''18'
$ws.Range('D41').Value = '''0.006987'
$ws.Range('G41').Value = '''18'
$ws.Range('G42').Value = '''18'
$ws.Range('D43').Value = '''0.002891'
$ws.Range('G43').Value = '''18'
$ws.Range('D44').Value = '''0.007473'
$ws.Range('G44').Value = '''18'
$ws.Range('D45').Value = '''0.00005939'
$ws.Range('G45').Value = '''18'
$ws.Range('D46').Value = '''0.00000000751'
$ws.Range('G46').Value = '''18'
$ws.Range('D47').Value = '''0.0005505'
$ws.Range('G47').Value = '''18'
$ws.Range('D48').Value = '''0.6831'
$ws.Range('G48').Value = '''18'
$ws.Range('D49').Value = '''0.04699'
$ws.Range('E49').Value = '48BOLOBOLOBestin24h'
$ws.Range('G49').Value = '''18'
$ws.Range('D50').Value = '''0.00002102'
$ws.Range('G50').Value = '''18'
$ws.Range('D51').Value = '''0.01011'
$ws.Range('G51').Value = '''18'
